# Updates Price (D) and Volume(1h) (E) columns for the cryptos worksheet
# to the latest scraped values (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.979.15"
$ws.Range("E2").Value = "  +2.66%  "
# Row 3
$ws.Range("D3").Value = "3.743.27"
$ws.Range("E3").Value = "  +2.23%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("D5").Value = "'602.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.42%  "
# Row 6
$ws.Range("D6").Value = "'168.86"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "
# Row 7
$ws.Range("D7").Value = "3.740.96"
$ws.Range("E7").Value = "  +2.23%  "
# Row 8
$ws.Range("E8").Value = "  -0.01%  "
# Row 9
$ws.Range("E9").Value = "  +2.11%  "
# Row 10
$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.04%  "
# Row 11
$ws.Range("D11").Value = "'6.36"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.48%  "
# Row 12
$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
# Row 13
$ws.Range("D13").Value = "'38.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.28%  "
# Row 14
$ws.Range("D14").Value = "'0.0000249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.15%  "
# Row 15
$ws.Range("D15").Value = "4.366.77"
$ws.Range("E15").Value = "  +1.95%  "
# Row 16
$ws.Range("D16").Value = "3.735.94"
$ws.Range("E16").Value = "  +1.93%  "
# Row 17
$ws.Range("D17").Value = "68.943.43"
$ws.Range("E17").Value = "  +2.58%  "
# Row 18
$ws.Range("E18").Value = "  +1.98%  "
# Row 19
$ws.Range("E19").Value = "  +0.28%  "
# Row 20
$ws.Range("D20").Value = "'17.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.91%  "
# Row 21
$ws.Range("D21").Value = "'10.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +19.15%  "
# Row 22
$ws.Range("D22").Value = "'494.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.59%  "
# Row 23
$ws.Range("D23").Value = "'0.729"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.77%  "
# Row 24
$ws.Range("D24").Value = "'0.0000156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.19%  "
# Row 25
$ws.Range("D25").Value = "'85.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.08%  "
# Row 26
$ws.Range("D26").Value = "'2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.47%  "
# Row 27
$ws.Range("E27").Value = "  +2.02%  "
# Row 28
$ws.Range("D28").Value = "'10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.67%  "
# Row 29
$ws.Range("E29").Value = "  +0.44%  "
# Row 30
$ws.Range("D30").Value = "'2.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.23%  "
# Row 31
$ws.Range("E31").Value = "  +2.10%  "
# Row 32
$ws.Range("E32").Value = "  +3.92%  "
# Row 33
$ws.Range("D33").Value = "'31.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.22%  "
# Row 34
$ws.Range("D34").Value = "3.886.58"
$ws.Range("E34").Value = "  +2.16%  "
# Row 35
$ws.Range("D35").Value = "'0.109"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.68%  "
# Row 36
$ws.Range("D36").Value = "3.675.40"
$ws.Range("E36").Value = "  +2.09%  "
# Row 37
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
# Row 38
$ws.Range("D38").Value = "'1.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.94%  "
# Row 39
$ws.Range("D39").Value = "'5.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "
# Row 40
$ws.Range("E40").Value = "  +1.78%  "
# Row 41
$ws.Range("E41").Value = "  +0.33%  "
# Row 42
$ws.Range("D42").Value = "'3.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.22%  "
# Row 43
$ws.Range("D43").Value = "'437.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.55%  "
# Row 44
$ws.Range("D44").Value = "'48.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.55%  "
# Row 45
$ws.Range("E45").Value = "  +3.17%  "
# Row 46
$ws.Range("E46").Value = "  +1.68%  "
# Row 48
$ws.Range("D48").Value = "'40.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.18%  "
# Row 49
$ws.Range("D49").Value = "'141.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.52%  "
# Row 50
$ws.Range("E50").Value = "  +2.63%  "
# Row 51
$ws.Range("D51").Value = "2.773.59"
$ws.Range("E51").Value = "  +0.98%  "
